$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(3.7060845695375324, 1.7903722316740134, -1.488474904563008, 2.6983668227749877)"
$ws.Range("C2").Value = "NIG(1.4311093228592093, 1.1077452835581276, 4.4311343483391425, 5.6466350598876724)"
$ws.Range("D2").Value = "NIG(0.6326850511674789, 0.4777104210397065, 0.7614956903478163, 2.9078787344582473)"
$ws.Range("E2").Value = "NIG(1.7695756864138863, 1.5617853773619146, 2.7951688189683335, 5.573084384464874)"
